$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.06469664608609
$ws.Cells.Item(2, 4).Value2 = 1.066648690080083
$ws.Cells.Item(2, 5).Value2 = 1.068170168690908
$ws.Cells.Item(2, 6).Value2 = 1.078445664253329
$ws.Cells.Item(2, 9).Value2 = 1.053419792298343
$ws.Cells.Item(2, 10).Value2 = 1.069655561025508
$ws.Cells.Item(2, 11).Value2 = 1.069358672877275
$ws.Cells.Item(2, 12).Value2 = 1.070876068790684
$ws.Cells.Item(2, 13).Value2 = 1.081124317237727
$ws.Cells.Item(2, 14).Value2 = 1.071174594613372
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.065815776577428
$ws.Cells.Item(3, 4).Value2 = 1.067520476643631
$ws.Cells.Item(3, 5).Value2 = 1.069177957709172
$ws.Cells.Item(3, 6).Value2 = 1.079458962644476
$ws.Cells.Item(3, 9).Value2 = 1.053752420509926
$ws.Cells.Item(3, 10).Value2 = 1.070429201011841
$ws.Cells.Item(3, 11).Value2 = 1.070045952391182
$ws.Cells.Item(3, 12).Value2 = 1.071699307232123
$ws.Cells.Item(3, 13).Value2 = 1.081955022508909
$ws.Cells.Item(3, 14).Value2 = 1.071949333257223
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.066540024701919
$ws.Cells.Item(4, 4).Value2 = 1.068084695376896
$ws.Cells.Item(4, 5).Value2 = 1.069830941647321
$ws.Cells.Item(4, 6).Value2 = 1.080115232197489
$ws.Cells.Item(4, 9).Value2 = 1.053966563001247
$ws.Cells.Item(4, 10).Value2 = 1.07092929778898
$ws.Cells.Item(4, 11).Value2 = 1.070490140182881
$ws.Cells.Item(4, 12).Value2 = 1.072232258031573
$ws.Cells.Item(4, 13).Value2 = 1.082492525160479
$ws.Cells.Item(4, 14).Value2 = 1.072450140229148
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.06684452197394
$ws.Cells.Item(5, 4).Value2 = 1.068321919978064
$ws.Cells.Item(5, 5).Value2 = 1.0701056655373
$ws.Cells.Item(5, 6).Value2 = 1.080391270611301
$ws.Cells.Item(5, 9).Value2 = 1.054056327421755
$ws.Cells.Item(5, 10).Value2 = 1.071139418730129
$ws.Cells.Item(5, 11).Value2 = 1.070676750154113
$ws.Cells.Item(5, 12).Value2 = 1.072456372506387
$ws.Cells.Item(5, 13).Value2 = 1.082718486378841
$ws.Cells.Item(5, 14).Value2 = 1.072660559566134
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.066895649787705
$ws.Cells.Item(6, 4).Value2 = 1.068361752607593
$ws.Cells.Item(6, 5).Value2 = 1.070151805153831
$ws.Cells.Item(6, 6).Value2 = 1.080437627031505
$ws.Cells.Item(6, 9).Value2 = 1.054071383950476
$ws.Cells.Item(6, 10).Value2 = 1.071174691937882
$ws.Cells.Item(6, 11).Value2 = 1.070708075362031
$ws.Cells.Item(6, 12).Value2 = 1.072494005934555
$ws.Cells.Item(6, 13).Value2 = 1.082756425965943
$ws.Cells.Item(6, 14).Value2 = 1.072695882865888
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.066544093317048
$ws.Cells.Item(7, 4).Value2 = 1.068087865077819
$ws.Cells.Item(7, 5).Value2 = 1.069834611699071
$ws.Cells.Item(7, 6).Value2 = 1.080118920076096
$ws.Cells.Item(7, 9).Value2 = 1.053967763463805
$ws.Cells.Item(7, 10).Value2 = 1.070932105905829
$ws.Cells.Item(7, 11).Value2 = 1.070492634172014
$ws.Cells.Item(7, 12).Value2 = 1.072235252417195
$ws.Cells.Item(7, 13).Value2 = 1.082495544484944
$ws.Cells.Item(7, 14).Value2 = 1.072452952333845
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.065074841872495
$ws.Cells.Item(8, 4).Value2 = 1.066943290491757
$ws.Cells.Item(8, 5).Value2 = 1.06851057379147
$ws.Cells.Item(8, 6).Value2 = 1.078787989179359
$ws.Cells.Item(8, 9).Value2 = 1.053532431241602
$ws.Cells.Item(8, 10).Value2 = 1.069917119792189
$ws.Cells.Item(8, 11).Value2 = 1.069591051432748
$ws.Cells.Item(8, 12).Value2 = 1.071154231901216
$ws.Cells.Item(8, 13).Value2 = 1.081405061770649
$ws.Cells.Item(8, 14).Value2 = 1.071436524823503
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.062486548631896
$ws.Cells.Item(9, 4).Value2 = 1.064927297153347
$ws.Cells.Item(9, 5).Value2 = 1.066184189984385
$ws.Cells.Item(9, 6).Value2 = 1.076447314422221
$ws.Cells.Item(9, 9).Value2 = 1.052756974227622
$ws.Cells.Item(9, 10).Value2 = 1.068124756240409
$ws.Cells.Item(9, 11).Value2 = 1.067998310702057
$ws.Cells.Item(9, 12).Value2 = 1.069251348027599
$ws.Cells.Item(9, 13).Value2 = 1.079483357670309
$ws.Cells.Item(9, 14).Value2 = 1.069641615909892
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.060761472020063
$ws.Cells.Item(10, 4).Value2 = 1.06358392174334
$ws.Cells.Item(10, 5).Value2 = 1.06463782801161
$ws.Cells.Item(10, 6).Value2 = 1.074889975696942
$ws.Cells.Item(10, 9).Value2 = 1.052234395154464
$ws.Cells.Item(10, 10).Value2 = 1.066927272486909
$ws.Cells.Item(10, 11).Value2 = 1.066933777063444
$ws.Cells.Item(10, 12).Value2 = 1.067984133383858
$ws.Cells.Item(10, 13).Value2 = 1.078202146540318
$ws.Cells.Item(10, 14).Value2 = 1.068442431592109
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.060014594757611
$ws.Cells.Item(11, 4).Value2 = 1.063002375839722
$ws.Cells.Item(11, 5).Value2 = 1.063969323318995
$ws.Cells.Item(11, 6).Value2 = 1.074216372478352
$ws.Cells.Item(11, 9).Value2 = 1.05200678291434
$ws.Cells.Item(11, 10).Value2 = 1.066408137137074
$ws.Cells.Item(11, 11).Value2 = 1.066472181216271
$ws.Cells.Item(11, 12).Value2 = 1.067435745415631
$ws.Cells.Item(11, 13).Value2 = 1.077647352663457
$ws.Cells.Item(11, 14).Value2 = 1.067922559010532
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.059737184415686
$ws.Cells.Item(12, 4).Value2 = 1.062786385694011
$ws.Cells.Item(12, 5).Value2 = 1.063721173565889
$ws.Cells.Item(12, 6).Value2 = 1.073966276845992
$ws.Cells.Item(12, 9).Value2 = 1.051922037430606
$ws.Cells.Item(12, 10).Value2 = 1.06621521427434
$ws.Cells.Item(12, 11).Value2 = 1.066300626931018
$ws.Cells.Item(12, 12).Value2 = 1.067232098675574
$ws.Cells.Item(12, 13).Value2 = 1.077441274472024
$ws.Cells.Item(12, 14).Value2 = 1.067729362175204
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.05979668929346
$ws.Cells.Item(13, 4).Value2 = 1.062832715319645
$ws.Cells.Item(13, 5).Value2 = 1.063774395164025
$ws.Cells.Item(13, 6).Value2 = 1.07401991820004
$ws.Cells.Item(13, 9).Value2 = 1.051940224682932
$ws.Cells.Item(13, 10).Value2 = 1.066256601089216
$ws.Cells.Item(13, 11).Value2 = 1.066337430298864
$ws.Cells.Item(13, 12).Value2 = 1.067275779372815
$ws.Cells.Item(13, 13).Value2 = 1.077485479079539
$ws.Cells.Item(13, 14).Value2 = 1.067770807764105
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.059991663668653
$ws.Cells.Item(14, 4).Value2 = 1.062984521579152
$ws.Cells.Item(14, 5).Value2 = 1.063948807862794
$ws.Cells.Item(14, 6).Value2 = 1.074195697238029
$ws.Cells.Item(14, 9).Value2 = 1.051999781909696
$ws.Cells.Item(14, 10).Value2 = 1.066392191968827
$ws.Cells.Item(14, 11).Value2 = 1.066458002465115
$ws.Cells.Item(14, 12).Value2 = 1.067418910912234
$ws.Cells.Item(14, 13).Value2 = 1.077630318237296
$ws.Cells.Item(14, 14).Value2 = 1.067906591198316
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.060111795533775
$ws.Cells.Item(15, 4).Value2 = 1.063078057332874
$ws.Cells.Item(15, 5).Value2 = 1.064056290847003
$ws.Cells.Item(15, 6).Value2 = 1.074304015166051
$ws.Cells.Item(15, 9).Value2 = 1.052036450558086
$ws.Cells.Item(15, 10).Value2 = 1.066475721671931
$ws.Cells.Item(15, 11).Value2 = 1.066532278097554
$ws.Cells.Item(15, 12).Value2 = 1.067507105473612
$ws.Cells.Item(15, 13).Value2 = 1.077719558022796
$ws.Cells.Item(15, 14).Value2 = 1.06799023952318
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.060811041685621
$ws.Cells.Item(16, 4).Value2 = 1.063622520090177
$ws.Cells.Item(16, 5).Value2 = 1.064682217256646
$ws.Cells.Item(16, 6).Value2 = 1.074934696021908
$ws.Cells.Item(16, 9).Value2 = 1.052249472961144
$ws.Cells.Item(16, 10).Value2 = 1.066961712775418
$ws.Cells.Item(16, 11).Value2 = 1.066964398052711
$ws.Cells.Item(16, 12).Value2 = 1.06802053494111
$ws.Cells.Item(16, 13).Value2 = 1.078238965949775
$ws.Cells.Item(16, 14).Value2 = 1.068476920789778
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.061249684510533
$ws.Cells.Item(17, 4).Value2 = 1.063964086014133
$ws.Cells.Item(17, 5).Value2 = 1.065075133823199
$ws.Cells.Item(17, 6).Value2 = 1.075330502316125
$ws.Cells.Item(17, 9).Value2 = 1.052382739631135
$ws.Cells.Item(17, 10).Value2 = 1.067266396840308
$ws.Cells.Item(17, 11).Value2 = 1.067235282584017
$ws.Cells.Item(17, 12).Value2 = 1.068342682684694
$ws.Cells.Item(17, 13).Value2 = 1.078564771275472
$ws.Cells.Item(17, 14).Value2 = 1.068782037540988
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.061505546403412
$ws.Cells.Item(18, 4).Value2 = 1.064163329569837
$ws.Cells.Item(18, 5).Value2 = 1.065304419612725
$ws.Cells.Item(18, 6).Value2 = 1.07556144030727
$ws.Cells.Item(18, 9).Value2 = 1.052460343303636
$ws.Cells.Item(18, 10).Value2 = 1.067444054382931
$ws.Cells.Item(18, 11).Value2 = 1.06739322264304
$ws.Cells.Item(18, 12).Value2 = 1.068530617290626
$ws.Cells.Item(18, 13).Value2 = 1.078754805975694
$ws.Cells.Item(18, 14).Value2 = 1.068959947377699
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.061592790238143
$ws.Cells.Item(19, 4).Value2 = 1.064231268830757
$ws.Cells.Item(19, 5).Value2 = 1.065382617805519
$ws.Cells.Item(19, 6).Value2 = 1.075640196209717
$ws.Cells.Item(19, 9).Value2 = 1.052486782352212
$ws.Cells.Item(19, 10).Value2 = 1.067504620924318
$ws.Cells.Item(19, 11).Value2 = 1.067447065554159
$ws.Cells.Item(19, 12).Value2 = 1.06859470343247
$ws.Cells.Item(19, 13).Value2 = 1.078819602549132
$ws.Cells.Item(19, 14).Value2 = 1.069020599930522
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.061202621365678
$ws.Cells.Item(20, 4).Value2 = 1.063927437780654
$ws.Cells.Item(20, 5).Value2 = 1.065032966813151
$ws.Cells.Item(20, 6).Value2 = 1.075288028713946
$ws.Cells.Item(20, 9).Value2 = 1.052368454672459
$ws.Cells.Item(20, 10).Value2 = 1.067233713281058
$ws.Cells.Item(20, 11).Value2 = 1.067206225681978
$ws.Cells.Item(20, 12).Value2 = 1.068308116043597
$ws.Cells.Item(20, 13).Value2 = 1.078529815656908
$ws.Cells.Item(20, 14).Value2 = 1.068749307567335
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.059934248218274
$ws.Cells.Item(21, 4).Value2 = 1.062939817813469
$ws.Cells.Item(21, 5).Value2 = 1.063897443177158
$ws.Cells.Item(21, 6).Value2 = 1.074143931637848
$ws.Cells.Item(21, 9).Value2 = 1.051982249316614
$ws.Cells.Item(21, 10).Value2 = 1.066352266399689
$ws.Cells.Item(21, 11).Value2 = 1.066422499646109
$ws.Cells.Item(21, 12).Value2 = 1.067376760887629
$ws.Cells.Item(21, 13).Value2 = 1.077587666798461
$ws.Cells.Item(21, 14).Value2 = 1.067866608930291
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.059136847271664
$ws.Cells.Item(22, 4).Value2 = 1.062318988718276
$ws.Cells.Item(22, 5).Value2 = 1.063184435648785
$ws.Cells.Item(22, 6).Value2 = 1.073425232372301
$ws.Cells.Item(22, 9).Value2 = 1.05173826858101
$ws.Cells.Item(22, 10).Value2 = 1.065797528300228
$ws.Cells.Item(22, 11).Value2 = 1.065929178559211
$ws.Cells.Item(22, 12).Value2 = 1.066791464371518
$ws.Cells.Item(22, 13).Value2 = 1.076995282866196
$ws.Cells.Item(22, 14).Value2 = 1.0673110830391
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.05955955768088
$ws.Cells.Item(23, 4).Value2 = 1.062648089889377
$ws.Cells.Item(23, 5).Value2 = 1.063562325089902
$ws.Cells.Item(23, 6).Value2 = 1.07380616771819
$ws.Cells.Item(23, 9).Value2 = 1.051867717171598
$ws.Cells.Item(23, 10).Value2 = 1.06609165646889
$ws.Cells.Item(23, 11).Value2 = 1.066190750623028
$ws.Cells.Item(23, 12).Value2 = 1.067101714181734
$ws.Cells.Item(23, 13).Value2 = 1.077309318422768
$ws.Cells.Item(23, 14).Value2 = 1.067605628903498
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.061223887147082
$ws.Cells.Item(24, 4).Value2 = 1.063943997496727
$ws.Cells.Item(24, 5).Value2 = 1.065052019947906
$ws.Cells.Item(24, 6).Value2 = 1.075307220487431
$ws.Cells.Item(24, 9).Value2 = 1.052374909827432
$ws.Cells.Item(24, 10).Value2 = 1.067248481758383
$ws.Cells.Item(24, 11).Value2 = 1.067219355437671
$ws.Cells.Item(24, 12).Value2 = 1.068323735124097
$ws.Cells.Item(24, 13).Value2 = 1.078545610602478
$ws.Cells.Item(24, 14).Value2 = 1.068764097017592
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.063155602668453
$ws.Cells.Item(25, 4).Value2 = 1.065448371879091
$ws.Cells.Item(25, 5).Value2 = 1.066784814247737
$ws.Cells.Item(25, 6).Value2 = 1.077051888163471
$ws.Cells.Item(25, 9).Value2 = 1.052958437129463
$ws.Cells.Item(25, 10).Value2 = 1.068588578943486
$ws.Cells.Item(25, 11).Value2 = 1.068410549955966
$ws.Cells.Item(25, 12).Value2 = 1.069743048421395
$ws.Cells.Item(25, 13).Value2 = 1.07998017858996
$ws.Cells.Item(25, 14).Value2 = 1.070106097294409
